# Apply NATMI LR-pair recalculation (Agrn-Dag1) per Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.59657533333333
$ws.Range("H2").Value = 46.789726
$ws.Range("I2").Value = 0.4757744772251148
$ws.Range("J2").Value = 0.475774477225115
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.779764
$ws.Range("N2").Value = 35.339292
$ws.Range("O2").Value = 0.1028447940505417
$ws.Range("P2").Value = 0.1028447940505417
$ws.Range("Q2").Value = 183.723976634888
$ws.Range("R2").Value = 1653.515789713992
$ws.Range("S2").Value = 0.04893092812472109
$ws.Range("T2").Value = 0.04893092812472109

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.59657533333333
$ws.Range("H3").Value = 46.789726
$ws.Range("I3").Value = 0.4757744772251148
$ws.Range("J3").Value = 0.475774477225115
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 43.841352
$ws.Range("N3").Value = 131.524056
$ws.Range("O3").Value = 0.3827627461243965
$ws.Range("P3").Value = 0.3827627461243964
$ws.Range("Q3").Value = 683.774949183184
$ws.Range("R3").Value = 6153.974542648656
$ws.Range("S3").Value = 0.1821087454385841
$ws.Range("T3").Value = 0.1821087454385841

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.59657533333333
$ws.Range("H4").Value = 46.789726
$ws.Range("I4").Value = 0.4757744772251148
$ws.Range("J4").Value = 0.475774477225115
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 58.91811866666666
$ws.Range("N4").Value = 176.754356
$ws.Range("O4").Value = 0.5143924598250619
$ws.Range("P4").Value = 0.5143924598250619
$ws.Range("Q4").Value = 918.9208762829395
$ws.Range("R4").Value = 8270.287886546457
$ws.Range("S4").Value = 0.2447348036618097
$ws.Range("T4").Value = 0.2447348036618098

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.399531333333333
$ws.Range("H5").Value = 19.198594
$ws.Range("I5").Value = 0.1952180917624358
$ws.Range("J5").Value = 0.1952180917624358
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.779764
$ws.Range("N5").Value = 35.339292
$ws.Range("O5").Value = 0.1028447940505417
$ws.Range("P5").Value = 0.1028447940505417
$ws.Range("Q5").Value = 75.384968817272
$ws.Range("R5").Value = 678.464719355448
$ws.Range("S5").Value = 0.02007716444224746
$ws.Range("T5").Value = 0.02007716444224746

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.399531333333333
$ws.Range("H6").Value = 19.198594
$ws.Range("I6").Value = 0.1952180917624358
$ws.Range("J6").Value = 0.1952180917624358
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 43.841352
$ws.Range("N6").Value = 131.524056
$ws.Range("O6").Value = 0.3827627461243965
$ws.Range("P6").Value = 0.3827627461243964
$ws.Range("Q6").Value = 280.564105819696
$ws.Range("R6").Value = 2525.076952377264
$ws.Range("S6").Value = 0.07472221289615434
$ws.Range("T6").Value = 0.07472221289615434

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.399531333333333
$ws.Range("H7").Value = 19.198594
$ws.Range("I7").Value = 0.1952180917624358
$ws.Range("J7").Value = 0.1952180917624358
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 58.91811866666666
$ws.Range("N7").Value = 176.754356
$ws.Range("O7").Value = 0.5143924598250619
$ws.Range("P7").Value = 0.5143924598250619
$ws.Range("Q7").Value = 377.0483465083848
$ws.Range("R7").Value = 3393.435118575464
$ws.Range("S7").Value = 0.100418714424034
$ws.Range("T7").Value = 0.100418714424034

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.78533933333333
$ws.Range("H8").Value = 32.356018
$ws.Range("I8").Value = 0.3290074310124493
$ws.Range("J8").Value = 0.3290074310124493
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 11.779764
$ws.Range("N8").Value = 35.339292
$ws.Range("O8").Value = 0.1028447940505417
$ws.Range("P8").Value = 0.1028447940505417
$ws.Range("Q8").Value = 127.048752006584
$ws.Range("R8").Value = 1143.438768059256
$ws.Range("S8").Value = 0.03383670148357316
$ws.Range("T8").Value = 0.03383670148357316

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.78533933333333
$ws.Range("H9").Value = 32.356018
$ws.Range("I9").Value = 0.3290074310124493
$ws.Range("J9").Value = 0.3290074310124493
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 43.841352
$ws.Range("N9").Value = 131.524056
$ws.Range("O9").Value = 0.3827627461243965
$ws.Range("P9").Value = 0.3827627461243964
$ws.Range("Q9").Value = 472.843858152112
$ws.Range("R9").Value = 4255.594723369008
$ws.Range("S9").Value = 0.125931787789658
$ws.Range("T9").Value = 0.125931787789658

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.78533933333333
$ws.Range("H10").Value = 32.356018
$ws.Range("I10").Value = 0.3290074310124493
$ws.Range("J10").Value = 0.3290074310124493
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 58.91811866666666
$ws.Range("N10").Value = 176.754356
$ws.Range("O10").Value = 0.5143924598250619
$ws.Range("P10").Value = 0.5143924598250619
$ws.Range("Q10").Value = 635.4519027016008
$ws.Range("R10").Value = 5719.067124314408
$ws.Range("S10").Value = 0.1692389417392181
$ws.Range("T10").Value = 0.1692389417392182

